# Generate Report for Handback
# Replace the two handed-back source files
#   9dbdf18d-2548-4aeb-873e-51fa422ec7d3.md  -> 786ce09c-c2f1-4abb-bc18-5632705c6788.md
#   ca8621fd-5cbb-4df7-a2b7-a492f2ab73fe.md  -> ffff3837dc51-0456-44ef-a254-0496531384a2.md
# and refresh the associated xliff correspondence file names / timestamps
# across the Overview, zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: Overview
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$ws1.Range("B2").Value = "e2e\786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$ws1.Range("G2").Value = "2016-08-21 07:07:52"

$ws1.Range("A3").Value = "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$ws1.Range("B3").Value = "e2e\ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$ws1.Range("G3").Value = "2016-08-21 07:07:52"

# Refresh hyperlinks on column B so their display text matches the new file names
# (the underlying link targets are left untouched).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd085d5ebb0a5b3b826883c25e63a49d7e3edcf2/e2e/9dbdf18d-2548-4aeb-873e-51fa422ec7d3.md", "", "", "e2e\786ce09c-c2f1-4abb-bc18-5632705c6788.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd085d5ebb0a5b3b826883c25e63a49d7e3edcf2/e2e/ca8621fd-5cbb-4df7-a2b7-a492f2ab73fe.md", "", "", "e2e\ffff3837dc51-0456-44ef-a254-0496531384a2.md")

# ----------------------------------------------------------------------
# Sheet 2: zh-cn
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$ws2.Range("G2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-21 07:07:48"
$ws2.Range("I2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$ws2.Range("J2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-21 07:08:09"

$ws2.Range("A3").Value = "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$ws2.Range("G3").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-21 07:07:48"
$ws2.Range("I3").Value = "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$ws2.Range("J3").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-21 07:08:09"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd085d5ebb0a5b3b826883c25e63a49d7e3edcf2/e2e/9dbdf18d-2548-4aeb-873e-51fa422ec7d3.md", "", "", "786ce09c-c2f1-4abb-bc18-5632705c6788.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2cd937bc94af97bdb42381fac907755f6c5d5b06/e2e/9dbdf18d-2548-4aeb-873e-51fa422ec7d3.md", "", "", "786ce09c-c2f1-4abb-bc18-5632705c6788.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd085d5ebb0a5b3b826883c25e63a49d7e3edcf2/e2e/ca8621fd-5cbb-4df7-a2b7-a492f2ab73fe.md", "", "", "ffff3837dc51-0456-44ef-a254-0496531384a2.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2cd937bc94af97bdb42381fac907755f6c5d5b06/e2e/ca8621fd-5cbb-4df7-a2b7-a492f2ab73fe.md", "", "", "ffff3837dc51-0456-44ef-a254-0496531384a2.md")

# ----------------------------------------------------------------------
# Sheet 3: de-de
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$ws3.Range("G2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-21 07:07:52"
$ws3.Range("I2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$ws3.Range("J2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-21 07:08:16"

$ws3.Range("A3").Value = "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$ws3.Range("G3").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-21 07:07:52"
$ws3.Range("I3").Value = "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$ws3.Range("J3").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-21 07:08:16"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd085d5ebb0a5b3b826883c25e63a49d7e3edcf2/e2e/9dbdf18d-2548-4aeb-873e-51fa422ec7d3.md", "", "", "786ce09c-c2f1-4abb-bc18-5632705c6788.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bc6281a420015df7c852307faa7889ebdf523bb1/e2e/9dbdf18d-2548-4aeb-873e-51fa422ec7d3.md", "", "", "786ce09c-c2f1-4abb-bc18-5632705c6788.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd085d5ebb0a5b3b826883c25e63a49d7e3edcf2/e2e/ca8621fd-5cbb-4df7-a2b7-a492f2ab73fe.md", "", "", "ffff3837dc51-0456-44ef-a254-0496531384a2.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bc6281a420015df7c852307faa7889ebdf523bb1/e2e/ca8621fd-5cbb-4df7-a2b7-a492f2ab73fe.md", "", "", "ffff3837dc51-0456-44ef-a254-0496531384a2.md")
